$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the password string value in B2 with a numeric value.
$ws.Range("B2").Value = 12345

# Move the active cell selection from E2 to E5 (leaves no trace of old value).
$ws.Range("E5").Select()
